$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the newTerm for "casual" (row 19): was "introduced", should be "uncertain"
$ws.Range("B19").Value = "uncertain"

# Fix origTerm text for row 26: "established.indoors" -> "invasive"
$ws.Range("A26").Value = "invasive"

# Fix origTerm text for row 27: "invasive" -> "established  invasive"
$ws.Range("A27").Value = "established  invasive"

# Add three new rows with origTerm/newTerm pairs
$ws.Range("A28").Value = "naturalised, invasive"
$ws.Range("B28").Value = "introduced"

$ws.Range("A29").Value = "naturalised, cryptogenic, invasive"
$ws.Range("B29").Value = "introduced"

$ws.Range("A30").Value = "naturalised, cryptogenic"
$ws.Range("B30").Value = "introduced"

# Update the selection / scroll position to match the edited cell
$ws.Range("B19").Select()

$wb.Save()
